$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F6").Value = "2/40"
$ws.Range("E10").Value = -27.93
$ws.Range("F10").Value = "11/40"
$ws.Range("G10").Value = -80.42
$ws.Range("G11").Value = -134.5
$ws.Range("F12").Value = "19/40"
$ws.Range("G12").Value = -116.74
$ws.Range("G13").Value = -200.47
$ws.Range("G14").Value = -2.91
$ws.Range("G15").Value = -32.42
$ws.Range("F16").Value = "16/40"
$ws.Range("G16").Value = -10.56
$ws.Range("G17").Value = -34.03
$ws.Range("F18").Value = "10/40"
$ws.Range("G18").Value = -4.1
$ws.Range("G19").Value = -45.17
$ws.Range("F20").Value = "3/40"
$ws.Range("G20").Value = 7.25
$ws.Range("G21").Value = -4.48
$ws.Range("E22").Value = 11.11
$ws.Range("F22").Value = "25/40"
$ws.Range("G22").Value = 6.63
$ws.Range("G23").Value = -3.72
$ws.Range("F24").Value = "29/40"
$ws.Range("G24").Value = 5.17
$ws.Range("G25").Value = -21.77
$ws.Range("E26").Value = 17.88
$ws.Range("G26").Value = -3.9
$ws.Range("G27").Value = -26.54
$ws.Range("E28").Value = 19.66
$ws.Range("F28").Value = "17/40"
$ws.Range("G28").Value = -6.88
$ws.Range("G29").Value = -22.4
$ws.Range("E30").Value = 46.07
$ws.Range("F30").Value = "6/40"
$ws.Range("G30").Value = 23.67
$ws.Range("G31").Value = -31.2
$ws.Range("G32").Value = -14.54
$ws.Range("G33").Value = -34.11
$ws.Range("F34").Value = "4/40"
$ws.Range("G34").Value = 15.89
$ws.Range("G35").Value = -11.46
$ws.Range("G36").Value = 13.24
$ws.Range("G37").Value = -8.93
$ws.Range("F38").Value = "26/40"
$ws.Range("G38").Value = 1.38
$ws.Range("G39").Value = -7.41
$ws.Range("E40").Value = 14.88
$ws.Range("F40").Value = "22/40"
$ws.Range("G40").Value = 7.47
$ws.Range("G41").Value = -25.14
$ws.Range("F42").Value = "27/40"
$ws.Range("G42").Value = -15.46
$ws.Range("G43").Value = -30.76
$ws.Range("F44").Value = "31/40"
$ws.Range("G44").Value = -22.86
$ws.Range("G45").Value = -31.4
$ws.Range("F46").Value = "23/40"
$ws.Range("G46").Value = -17.4
$ws.Range("G47").Value = -30.27
$ws.Range("G48").Value = -6.78
$ws.Range("G49").Value = -33.41
$ws.Range("F50").Value = "36/40"
$ws.Range("G50").Value = -28.26
$ws.Range("G51").Value = -39.45
$ws.Range("F52").Value = "32/40"
$ws.Range("G52").Value = -32.36
$ws.Range("G53").Value = -38.98
$ws.Range("F54").Value = "28/40"
$ws.Range("G54").Value = -29.96
$ws.Range("G55").Value = -43.06
$ws.Range("E56").Value = 17.56
$ws.Range("F56").Value = "20/40"
$ws.Range("G56").Value = -25.5
$ws.Range("G57").Value = -37.76
$ws.Range("G58").Value = -11.14
$ws.Range("G59").Value = -35.57
$ws.Range("F60").Value = "15/40"
$ws.Range("G60").Value = -12.55
$ws.Range("G61").Value = -34.19
$ws.Range("G62").Value = -30.46
$ws.Range("G63").Value = -54.92
$ws.Range("G64").Value = -48.25
$ws.Range("G65").Value = -54.5
$ws.Range("F66").Value = "30/40"
$ws.Range("G66").Value = -45.93
$ws.Range("G67").Value = -51.19
$ws.Range("F68").Value = "24/40"
$ws.Range("G68").Value = -38.23
$ws.Range("G69").Value = -48.07
$ws.Range("F70").Value = "9/40"
$ws.Range("G70").Value = -18.07
$ws.Range("G71").Value = -46.04
$ws.Range("F72").Value = "5/40"
$ws.Range("G72").Value = 2.5
$ws.Range("G73").Value = -13.18
$ws.Range("E74").Value = 5.44
$ws.Range("F74").Value = "35/40"
$ws.Range("G74").Value = -7.74
$ws.Range("G75").Value = -13.47
$ws.Range("G76").Value = -10.09
$ws.Range("G77").Value = -29.7
$ws.Range("G78").Value = -23.18
$ws.Range("G79").Value = -48.35
$ws.Range("G80").Value = -9.26
$ws.Range("G81").Value = -12.53

Write-Host "Applied 106 cell updates to grov_rise_events sheet"
